$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the current row 44 (i.e. before the
# current row 45), shifting all subsequent rows down by two. This mirrors
# the diff, where the former rows 45-55 become rows 47-57 and two brand
# new rows of data are introduced at positions 45 and 46.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# New row 45: Alcachofa, Argentina(o), Primera
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44769
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112013
$ws.Cells.Item(45, 7).Value = "Alcachofa"
$ws.Cells.Item(45, 8).Value = "Argentina(o)"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 13000
$ws.Cells.Item(45, 12).Value = 14000
$ws.Cells.Item(45, 13).Value = 13500
$ws.Cells.Item(45, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(45, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(45, 16).Value = 270
$ws.Cells.Item(45, 17).Value = 50
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# New row 46: Alcachofa, Española, Primera
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 44769
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112013
$ws.Cells.Item(46, 7).Value = "Alcachofa"
$ws.Cells.Item(46, 8).Value = "Española"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 16000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 13).Value = 17000
$ws.Cells.Item(46, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 567
$ws.Cells.Item(46, 17).Value = 30
$ws.Cells.Item(46, 18).Value = "Hortaliza"
